# Updated cryptos list values (price & volume/1h) to reflect refreshed market data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.877.42'
$ws.Range("E2").Value = '  +3.20%  '
$ws.Range("D3").Value = '2.683.67'
$ws.Range("E3").Value = '  +1.71%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").Value = "'521.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.25%  '
$ws.Range("D6").Value = "'146.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.43%  '
$ws.Range("D7").Value = "'0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.21%  '
$ws.Range("D8").Value = "'0.579"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.66%  '
$ws.Range("D9").Value = '2.703.48'
$ws.Range("E9").Value = '  +1.45%  '
$ws.Range("D10").Value = "'6.43"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.12%  '
$ws.Range("D11").Value = "'0.106"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.49%  '
$ws.Range("D12").Value = "'0.341"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.89%  '
$ws.Range("E13").Value = '  +1.48%  '
$ws.Range("D14").Value = '3.158.21'
$ws.Range("E14").Value = '  +1.91%  '
$ws.Range("D15").Value = '60.638.22'
$ws.Range("E15").Value = '  +2.81%  '
$ws.Range("D16").Value = "'21.40"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.12%  '
$ws.Range("D17").Value = '2.755.77'
$ws.Range("E17").Value = '  +3.58%  '
$ws.Range("D18").Value = "'0.0000139"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.99%  '
$ws.Range("D19").Value = "'352.97"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.26%  '
$ws.Range("D20").Value = "'4.57"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.02%  '
$ws.Range("E21").Value = '  +1.38%  '
$ws.Range("E22").Value = '  +3.67%  '
$ws.Range("D23").Value = "'0.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.06%  '
$ws.Range("D24").Value = "'62.99"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.47%  '
$ws.Range("D25").Value = "'0.423"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.06%  '
$ws.Range("D26").Value = "'0.169"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.53%  '
$ws.Range("D27").Value = "'0.995"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.08%  '
$ws.Range("D28").Value = '0.0₃0820'
$ws.Range("E28").Value = '  +0.96%  '
$ws.Range("D29").Value = "'7.31"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.17%  '
$ws.Range("D30").Value = "'6.84"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.34%  '
$ws.Range("E31").Value = '  +0.09%  '
$ws.Range("D32").Value = "'19.12"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.71%  '
$ws.Range("D33").Value = "'1.59"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.76%  '
$ws.Range("D34").Value = "'149.40"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.63%  '
$ws.Range("D35").Value = "'4.31"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +6.54%  '
$ws.Range("E36").Value = '  +6.39%  '
$ws.Range("E37").Value = '  -8.81%  '
$ws.Range("E38").Value = '  +10.86%  '
$ws.Range("D39").Value = "'0.876"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.60%  '
$ws.Range("E40").Value = '  +0.30%  '
$ws.Range("D41").Value = "'3.72"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.28%  '
$ws.Range("D42").Value = "'284.37"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.17%  '
$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").Value = "'20.08"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.42%  '
$ws.Range("B44").Value = 'Stellar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D44").Value = "'0.0992"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.78%  '
$ws.Range("D45").Value = "'0.612"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.56%  '
$ws.Range("D46").Value = "'0.998"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.46%  '
$ws.Range("D47").Value = '2.133.79'
$ws.Range("E47").Value = '  +7.57%  '
$ws.Range("E48").Value = '  +0.49%  '
$ws.Range("D49").Value = "'4.87"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.72%  '
$ws.Range("E50").Value = '  +2.81%  '
$ws.Range("B51").Value = 'WhiteBITCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D51").Value = "'10.46"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.01%  '
